$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new taxonomy rows sourced from the CA Water Boards microplastics
#     methodology doc (rows 12-13) and a watres.2018 DOI reference (row 14) ---
$ws.Range("C12").Value = "Rubbery fragment"
$ws.Range("C13").Value = "Fiber bundle"

# New alias for the existing "Bead" row
$ws.Range("D9").Value = "Beads"

$ws.Range("C14").Value = "Other"
$ws.Range("A14").Value = "doi.org/10.1016/j.watres.2018.10.013"

$ws.Range("A12").Value = "https://www.waterboards.ca.gov/drinking_water/certlic/drinkingwater/documents/microplastics/swb-mp2-rev1.pdf"
$ws.Range("A13").Value = "https://www.waterboards.ca.gov/drinking_water/certlic/drinkingwater/documents/microplastics/swb-mp2-rev1.pdf"

# New aliases for the "Foam" row (Styrofoam / Polystyrene)
$ws.Range("E6").Value = "Styrofoam"
$ws.Range("F6").Value = "Polystyrene"

# Insert a new "Alias_4" column before the existing "Locked (Y/N)" column,
# shifting "Locked (Y/N)" to H and "Primary Alias" to I
$ws.Columns("G").Insert()
$ws.Range("G1").Value = "Alias_4"
$ws.Range("G6").Value = "Extruded polystyrene"

# Widen column D (Alias_1) to fit its longest entry
$ws.Columns("D").ColumnWidth = 15.1666666

# Match the final selection state
$ws.Range("G7").Select() | Out-Null
